# Update crypto price ("D") and volume-change ("E") columns for rows 2-51
# to reflect the refreshed feed values described in the commit
# "Updated cryptos list on Wed Apr 12 20:27:13 UTC 2023 with GitHub Actions".
#
# Values in column D frequently look like plain numbers (e.g. "1.099",
# "24.04") but must stay as literal text, exactly matching the original
# workbook's inline-string formatting (no numeric coercion / no rounding).
# We force text by temporarily applying a "@" (Text) number format before
# writing the value, then restore the cell to the "Normal" style so the
# saved file keeps the same (default) style index as every other data cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "29.904.91"; E = "  -0.78%  " },
    @{ Row = 3; D = "1.909.66"; E = "  +0.80%  " },
    @{ Row = 4; D = $null; E = "  -0.18%  " },
    @{ Row = 5; D = "319.90"; E = "  -1.80%  " },
    @{ Row = 6; D = $null; E = "  -0.11%  " },
    @{ Row = 7; D = "0.5048"; E = "  -2.41%  " },
    @{ Row = 8; D = "0.4048"; E = "  +1.05%  " },
    @{ Row = 9; D = $null; E = "  -1.90%  " },
    @{ Row = 10; D = "41.94"; E = "  -1.75%  " },
    @{ Row = 11; D = "1.099"; E = "  -1.68%  " },
    @{ Row = 12; D = "24.04"; E = "  +3.75%  " },
    @{ Row = 13; D = "1.906.96"; E = "  +0.33%  " },
    @{ Row = 14; D = "6.381"; E = "  -0.80%  " },
    @{ Row = 15; D = "7.212"; E = "  -1.63%  " },
    @{ Row = 16; D = "1.001"; E = "  -0.29%  " },
    @{ Row = 17; D = "92.09"; E = "  -2.70%  " },
    @{ Row = 18; D = $null; E = "  -1.45%  " },
    @{ Row = 19; D = "0.06512"; E = "  -2.21%  " },
    @{ Row = 20; D = "18.12"; E = "  -0.79%  " },
    @{ Row = 21; D = $null; E = "  -0.10%  " },
    @{ Row = 22; D = "5.942"; E = "  -0.08%  " },
    @{ Row = 23; D = "29.939.88"; E = "  -0.72%  " },
    @{ Row = 24; D = "11.28"; E = "  -0.01%  " },
    @{ Row = 25; D = "2.187"; E = "  -1.30%  " },
    @{ Row = 26; D = "22.08"; E = "  +0.74%  " },
    @{ Row = 27; D = "2.130.46"; E = "  +0.65%  " },
    @{ Row = 28; D = "161.76"; E = "  +0.34%  " },
    @{ Row = 29; D = "2.292"; E = "  -3.79%  " },
    @{ Row = 30; D = "128.66"; E = "  -0.25%  " },
    @{ Row = 31; D = $null; E = "  +3.03%  " },
    @{ Row = 32; D = $null; E = "  -2.04%  " },
    @{ Row = 33; D = "5.921"; E = "  -2.67%  " },
    @{ Row = 34; D = "3.798"; E = "  +2.63%  " },
    @{ Row = 35; D = "5.401"; E = "  +2.98%  " },
    @{ Row = 36; D = $null; E = "  -2.37%  " },
    @{ Row = 37; D = "0.06374"; E = "  -2.73%  " },
    @{ Row = 38; D = "0.2146"; E = "  -2.79%  " },
    @{ Row = 39; D = "1.194"; E = "  -1.69%  " },
    @{ Row = 40; D = "8.693"; E = "  -0.82%  " },
    @{ Row = 41; D = "0.6450"; E = "  -0.75%  " },
    @{ Row = 42; D = $null; E = "  -3.73%  " },
    @{ Row = 43; D = $null; E = "  -2.25%  " },
    @{ Row = 44; D = "2.208"; E = "  +7.49%  " },
    @{ Row = 45; D = "13.23"; E = "  -0.31%  " },
    @{ Row = 46; D = "0.6034"; E = "  -1.10%  " },
    @{ Row = 47; D = "3.634"; E = "  -2.00%  " },
    @{ Row = 48; D = "122.15"; E = "  -1.86%  " },
    @{ Row = 49; D = "1.205"; E = "  -2.73%  " },
    @{ Row = 50; D = "78.61"; E = "  -0.63%  " },
    @{ Row = 51; D = "1.123"; E = "  -3.43%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}

Write-Host "Applied $($updates.Count) row updates to cryptos sheet"
